# Hortaliza / Zanahoria - Terminal La Palmera de La Serena
# Insert a new weekly price observation row at row 426 (pushing all
# subsequent rows down by one, growing the used range from A1:R504 to
# A1:R505), then populate the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 426:504 down to 427:505, leaving a blank row 426 behind.
$ws.Rows.Item(426).Insert()

# Fill in the new observation row.
$ws.Range("A426").Value = 8
$ws.Range("B426").Value = "Terminal La Palmera de La Serena"
$ws.Range("C426").Value = "Coquimbo"
$ws.Range("D426").Value = 45015
$ws.Range("E426").Value = 4
$ws.Range("F426").Value = 100114013
$ws.Range("G426").Value = "Zanahoria"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 480
$ws.Range("K426").Value = 5500
$ws.Range("L426").Value = 6000
$ws.Range("M426").Value = 5750
$ws.Range("N426").Value = "`$/saco 20 kilos"
$ws.Range("O426").Value = "Provincia del Elquí"
$ws.Range("P426").Value = 288
$ws.Range("Q426").Value = 20
$ws.Range("R426").Value = "Hortaliza"
